# Auto-generated: apply numeric cell updates per the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 21000050
$ws.Cells.Item(6, 9).Value = 200059.4
$ws.Cells.Item(6, 11).Value = 600178.2
$ws.Cells.Item(6, 13).Value = -600066.2
$ws.Cells.Item(106, 8).Value = 3528.3333
$ws.Cells.Item(106, 9).Value = 3234.4
$ws.Cells.Item(106, 11).Value = 3234.4
$ws.Cells.Item(106, 13).Value = -2603.4
$ws.Cells.Item(113, 8).Value = 25652146
$ws.Cells.Item(113, 9).Value = 29415256
$ws.Cells.Item(113, 11).Value = 29415256
$ws.Cells.Item(113, 13).Value = -29412002
$ws.Cells.Item(132, 8).Value = 1717.8334
$ws.Cells.Item(132, 9).Value = 1717.8334
$ws.Cells.Item(132, 11).Value = 5153.5002
$ws.Cells.Item(132, 13).Value = -2623.5002
$ws.Cells.Item(135, 8).Value = 756.9655
$ws.Cells.Item(135, 9).Value = 779.1429000000001
$ws.Cells.Item(135, 11).Value = 7012.2861
$ws.Cells.Item(135, 13).Value = -4477.2861

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 8436.571
$ws.Cells.Item(61, 9).Value = 9201
$ws.Cells.Item(61, 10).Value = 3850
$ws.Cells.Item(61, 11).Value = 9201
$ws.Cells.Item(61, 12).Value = 3850
$ws.Cells.Item(61, 13).Value = -8989
$ws.Cells.Item(61, 14).Value = -4274
$ws.Cells.Item(74, 8).Value = 6598.75
$ws.Cells.Item(74, 9).Value = 3798.3333
$ws.Cells.Item(74, 10).Value = 15000
$ws.Cells.Item(74, 11).Value = 3798.3333
$ws.Cells.Item(74, 12).Value = 15000
$ws.Cells.Item(74, 13).Value = -2924.3333
$ws.Cells.Item(74, 14).Value = -16748
$ws.Cells.Item(77, 8).Value = 6598.75
$ws.Cells.Item(77, 9).Value = 3798.3333
$ws.Cells.Item(77, 10).Value = 15000
$ws.Cells.Item(77, 11).Value = 18991.6665
$ws.Cells.Item(77, 12).Value = 75000
$ws.Cells.Item(77, 13).Value = -14623.6665
$ws.Cells.Item(77, 14).Value = -83736
$ws.Cells.Item(102, 8).Value = 4118.9375
$ws.Cells.Item(102, 9).Value = 2993.0715
$ws.Cells.Item(102, 11).Value = 2993.0715
$ws.Cells.Item(102, 13).Value = -1371.0715
$ws.Cells.Item(132, 8).Value = 1744.742
$ws.Cells.Item(132, 9).Value = 1713.0769
$ws.Cells.Item(132, 10).Value = 1909.4
$ws.Cells.Item(132, 11).Value = 5139.2307
$ws.Cells.Item(132, 12).Value = 5728.200000000001
$ws.Cells.Item(132, 13).Value = -2609.2307
$ws.Cells.Item(132, 14).Value = -10788.2
$ws.Cells.Item(134, 8).Value = 173333.33
$ws.Cells.Item(134, 10).Value = 173333.33
$ws.Cells.Item(134, 12).Value = 173333.33
$ws.Cells.Item(134, 14).Value = -183473.33
$ws.Cells.Item(136, 8).Value = 8436.571
$ws.Cells.Item(136, 9).Value = 9201
$ws.Cells.Item(136, 10).Value = 3850
$ws.Cells.Item(136, 11).Value = 27603
$ws.Cells.Item(136, 12).Value = 11550
$ws.Cells.Item(136, 13).Value = -25053
$ws.Cells.Item(136, 14).Value = -16650

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(103, 8).Value = 24999
$ws.Cells.Item(103, 10).Value = 24999
$ws.Cells.Item(103, 12).Value = 24999
$ws.Cells.Item(103, 14).Value = -27343
$ws.Cells.Item(105, 8).Value = 3553.5557
$ws.Cells.Item(105, 9).Value = 3560.25
$ws.Cells.Item(105, 11).Value = 3560.25
$ws.Cells.Item(105, 13).Value = -1813.25
$ws.Cells.Item(107, 8).Value = 3333.7646
$ws.Cells.Item(107, 9).Value = 3090.3572
$ws.Cells.Item(107, 11).Value = 3090.3572
$ws.Cells.Item(107, 13).Value = -1170.3572
$ws.Cells.Item(134, 8).Value = 4132.968
$ws.Cells.Item(134, 9).Value = 3522.6072
$ws.Cells.Item(134, 10).Value = 9829.666999999999
$ws.Cells.Item(134, 11).Value = 10567.8216
$ws.Cells.Item(134, 12).Value = 29489.001
$ws.Cells.Item(134, 13).Value = -8032.821599999999
$ws.Cells.Item(134, 14).Value = -34559.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 10078.263
$ws.Cells.Item(58, 9).Value = 4848
$ws.Cells.Item(58, 10).Value = 12492.23
$ws.Cells.Item(58, 11).Value = 4848
$ws.Cells.Item(58, 12).Value = 12492.23
$ws.Cells.Item(58, 13).Value = -4645
$ws.Cells.Item(58, 14).Value = -12898.23
$ws.Cells.Item(105, 8).Value = 803.5833
$ws.Cells.Item(105, 9).Value = 860.55554
$ws.Cells.Item(105, 11).Value = 860.55554
$ws.Cells.Item(105, 13).Value = 886.44446
$ws.Cells.Item(136, 8).Value = 10078.263
$ws.Cells.Item(136, 9).Value = 4848
$ws.Cells.Item(136, 10).Value = 12492.23
$ws.Cells.Item(136, 11).Value = 14544
$ws.Cells.Item(136, 12).Value = 37476.69
$ws.Cells.Item(136, 13).Value = -11994
$ws.Cells.Item(136, 14).Value = -42576.69

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(15, 8).Value = 766.6667
$ws.Cells.Item(15, 9).Value = 766.6667
$ws.Cells.Item(15, 11).Value = 2300.0001
$ws.Cells.Item(15, 13).Value = -2160.0001
$ws.Cells.Item(31, 8).Value = 1273
$ws.Cells.Item(31, 9).Value = 1546
$ws.Cells.Item(31, 10).Value = 1000
$ws.Cells.Item(31, 11).Value = 4638
$ws.Cells.Item(31, 12).Value = 3000
$ws.Cells.Item(31, 13).Value = -4350
$ws.Cells.Item(31, 14).Value = -3576
$ws.Cells.Item(122, 8).Value = 8164.75
$ws.Cells.Item(122, 10).Value = 13541.667
$ws.Cells.Item(122, 12).Value = 121875.003
$ws.Cells.Item(122, 14).Value = -126775.003

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(7, 8).Value = 336063.34
$ws.Cells.Item(7, 10).Value = 336063.34
$ws.Cells.Item(7, 12).Value = 336063.34
$ws.Cells.Item(7, 14).Value = -336287.34
$ws.Cells.Item(8, 8).Value = 336063.34
$ws.Cells.Item(8, 10).Value = 336063.34
$ws.Cells.Item(8, 12).Value = 336063.34
$ws.Cells.Item(8, 14).Value = -336341.34
$ws.Cells.Item(107, 8).Value = 1376.75
$ws.Cells.Item(107, 9).Value = 561
$ws.Cells.Item(107, 10).Value = 1959.4286
$ws.Cells.Item(107, 11).Value = 561
$ws.Cells.Item(107, 12).Value = 1959.4286
$ws.Cells.Item(107, 13).Value = 1359
$ws.Cells.Item(107, 14).Value = -5799.4286
$ws.Cells.Item(132, 8).Value = 1631.7576
$ws.Cells.Item(132, 9).Value = 1369.8438
$ws.Cells.Item(132, 10).Value = 10013
$ws.Cells.Item(132, 11).Value = 4109.5314
$ws.Cells.Item(132, 12).Value = 30039
$ws.Cells.Item(132, 13).Value = -1579.5314
$ws.Cells.Item(132, 14).Value = -35099
$ws.Cells.Item(136, 8).Value = 39346.5
$ws.Cells.Item(136, 10).Value = 39346.5
$ws.Cells.Item(136, 12).Value = 118039.5
$ws.Cells.Item(136, 14).Value = -123139.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1623
$ws.Cells.Item(22, 9).Value = 1743.5
$ws.Cells.Item(22, 11).Value = 1743.5
$ws.Cells.Item(22, 13).Value = -1448.5
$ws.Cells.Item(27, 8).Value = 1623
$ws.Cells.Item(27, 9).Value = 1743.5
$ws.Cells.Item(27, 11).Value = 1743.5
$ws.Cells.Item(27, 13).Value = -1636.5
$ws.Cells.Item(82, 8).Value = 2399.3333
$ws.Cells.Item(82, 9).Value = 2099
$ws.Cells.Item(82, 11).Value = 2099
$ws.Cells.Item(82, 13).Value = -1738
$ws.Cells.Item(85, 8).Value = 2399.3333
$ws.Cells.Item(85, 9).Value = 2099
$ws.Cells.Item(85, 11).Value = 2099
$ws.Cells.Item(85, 13).Value = -851
$ws.Cells.Item(132, 8).Value = 8939.311
$ws.Cells.Item(132, 9).Value = 8871.813
$ws.Cells.Item(132, 11).Value = 26615.439
$ws.Cells.Item(132, 13).Value = -24085.439

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 3945.6924
$ws.Cells.Item(81, 9).Value = 2310.5557
$ws.Cells.Item(81, 10).Value = 7624.75
$ws.Cells.Item(81, 11).Value = 4621.1114
$ws.Cells.Item(81, 12).Value = 15249.5
$ws.Cells.Item(81, 13).Value = -3560.1114
$ws.Cells.Item(81, 14).Value = -17371.5
$ws.Cells.Item(84, 8).Value = 3945.6924
$ws.Cells.Item(84, 9).Value = 2310.5557
$ws.Cells.Item(84, 10).Value = 7624.75
$ws.Cells.Item(84, 11).Value = 23105.557
$ws.Cells.Item(84, 12).Value = 76247.5
$ws.Cells.Item(84, 13).Value = -17801.557
$ws.Cells.Item(84, 14).Value = -86855.5
$ws.Cells.Item(95, 8).Value = 64994.5
$ws.Cells.Item(95, 10).Value = 64994.5
$ws.Cells.Item(95, 12).Value = 64994.5
$ws.Cells.Item(95, 14).Value = -70486.5
